$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 275.83334
$ws.Range("I6").Value = 74.5
$ws.Range("J6").Value = 301
$ws.Range("K6").Value = 223.5
$ws.Range("L6").Value = 903
$ws.Range("M6").Value = -111.5
$ws.Range("N6").Value = -1127

$ws.Range("H15").Value = 954.14636
$ws.Range("I15").Value = 954.14636
$ws.Range("K15").Value = 2862.43908
$ws.Range("M15").Value = -2693.43908

$ws.Range("H51").Value = 73221000
$ws.Range("J51").Value = 8855.625
$ws.Range("L51").Value = 8855.625
$ws.Range("N51").Value = -9823.625

$ws.Range("H69").Value = 15427
$ws.Range("J69").Value = 9971
$ws.Range("L69").Value = 29913
$ws.Range("N69").Value = -31661

$ws.Range("H72").Value = 15427
$ws.Range("J72").Value = 9971
$ws.Range("L72").Value = 89739
$ws.Range("N72").Value = -98475

$ws.Range("H88").Value = 848.6
$ws.Range("I88").Value = 822
$ws.Range("J88").Value = 866.3333
$ws.Range("K88").Value = 822
$ws.Range("L88").Value = 866.3333
$ws.Range("M88").Value = -416
$ws.Range("N88").Value = -1678.3333

$ws.Range("H91").Value = 848.6
$ws.Range("I91").Value = 822
$ws.Range("J91").Value = 866.3333
$ws.Range("K91").Value = 822
$ws.Range("L91").Value = 866.3333
$ws.Range("M91").Value = 582
$ws.Range("N91").Value = -3674.3333

$ws.Range("H138").Value = 6054.95
$ws.Range("I138").Value = 3011.5
$ws.Range("J138").Value = 10620.125
$ws.Range("K138").Value = 9034.5
$ws.Range("L138").Value = 31860.375
$ws.Range("M138").Value = -3894.5
$ws.Range("N138").Value = -42140.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 2357833
$ws.Range("I6").Value = 2884356.8
$ws.Range("J6").Value = 515000.5
$ws.Range("K6").Value = 2884356.8
$ws.Range("L6").Value = 515000.5
$ws.Range("M6").Value = -2884183.8
$ws.Range("N6").Value = -515346.5

$ws.Range("H32").Value = 2630.581
$ws.Range("I32").Value = 2596.1694
$ws.Range("K32").Value = 2596.1694
$ws.Range("M32").Value = -2309.1694

$ws.Range("H45").Value = 4147.75
$ws.Range("I45").Value = 2818.5557
$ws.Range("J45").Value = 5856.7144
$ws.Range("K45").Value = 2818.5557
$ws.Range("L45").Value = 5856.7144
$ws.Range("M45").Value = -2441.5557
$ws.Range("N45").Value = -6610.7144

$ws.Range("H61").Value = 4086.973
$ws.Range("J61").Value = 10093.25
$ws.Range("L61").Value = 10093.25
$ws.Range("N61").Value = -10517.25

$ws.Range("H88").Value = 3994.4
$ws.Range("I88").Value = 3322.6667
$ws.Range("J88").Value = 5002
$ws.Range("K88").Value = 3322.6667
$ws.Range("L88").Value = 5002
$ws.Range("M88").Value = -2916.6667
$ws.Range("N88").Value = -5814

$ws.Range("H91").Value = 3994.4
$ws.Range("I91").Value = 3322.6667
$ws.Range("J91").Value = 5002
$ws.Range("K91").Value = 3322.6667
$ws.Range("L91").Value = 5002
$ws.Range("M91").Value = -1918.6667
$ws.Range("N91").Value = -7810

$ws.Range("H110").Value = 5497.5
$ws.Range("I110").Value = 4259.8335
$ws.Range("K110").Value = 4259.8335
$ws.Range("M110").Value = -2214.8335

$ws.Range("H132").Value = 3981.8857
$ws.Range("I132").Value = 2127.32
$ws.Range("K132").Value = 6381.960000000001
$ws.Range("M132").Value = -3851.960000000001

$ws.Range("H136").Value = 4086.973
$ws.Range("J136").Value = 10093.25
$ws.Range("L136").Value = 30279.75
$ws.Range("N136").Value = -35379.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H140").Value = 87000
$ws.Range("J140").Value = 87000
$ws.Range("L140").Value = 87000
$ws.Range("N140").Value = -97360

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 505000000
$ws.Range("J4").Value = 505000000
$ws.Range("L4").Value = 505000000
$ws.Range("N4").Value = -505000224

$ws.Range("H31").Value = 1768.7894
$ws.Range("I31").Value = 1286.2858
$ws.Range("J31").Value = 3119.8
$ws.Range("K31").Value = 1286.2858
$ws.Range("L31").Value = 3119.8
$ws.Range("M31").Value = -991.2858000000001
$ws.Range("N31").Value = -3709.8

$ws.Range("H34").Value = 1768.7894
$ws.Range("I34").Value = 1286.2858
$ws.Range("J34").Value = 3119.8
$ws.Range("K34").Value = 1286.2858
$ws.Range("L34").Value = 3119.8
$ws.Range("M34").Value = -1084.2858
$ws.Range("N34").Value = -3523.8

$ws.Range("H56").Value = 11379.8
$ws.Range("I56").Value = 7974.75
$ws.Range("K56").Value = 7974.75
$ws.Range("M56").Value = -7129.75

$ws.Range("H109").Value = 42428.332
$ws.Range("J109").Value = 42428.332
$ws.Range("L109").Value = 42428.332
$ws.Range("N109").Value = -44508.332

$ws.Range("H132").Value = 1952.6666
$ws.Range("I132").Value = 1649.8462
$ws.Range("K132").Value = 4949.5386
$ws.Range("M132").Value = -2419.5386

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 517.625
$ws.Range("I86").Value = 561
$ws.Range("J86").Value = 474.25
$ws.Range("K86").Value = 1683
$ws.Range("L86").Value = 1422.75
$ws.Range("M86").Value = -497
$ws.Range("N86").Value = -3794.75

$ws.Range("H89").Value = 517.625
$ws.Range("I89").Value = 561
$ws.Range("J89").Value = 474.25
$ws.Range("K89").Value = 5049
$ws.Range("L89").Value = 4268.25
$ws.Range("M89").Value = 879
$ws.Range("N89").Value = -16124.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 9513201
$ws.Range("J3").Value = 23775000
$ws.Range("L3").Value = 23775000
$ws.Range("N3").Value = -23775232

$ws.Range("H29").Value = 19999
$ws.Range("J29").Value = 19998
$ws.Range("L29").Value = 19998
$ws.Range("N29").Value = -20578

$ws.Range("H42").Value = 92446
$ws.Range("J42").Value = 92446
$ws.Range("L42").Value = 92446
$ws.Range("N42").Value = -93416

$ws.Range("H115").Value = 92446
$ws.Range("J115").Value = 92446
$ws.Range("L115").Value = 92446
$ws.Range("N115").Value = -94796

$ws.Range("H132").Value = 3853.6165
$ws.Range("I132").Value = 3374.3936
$ws.Range("K132").Value = 10123.1808
$ws.Range("M132").Value = -7593.1808

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()

$ws.Range("H40").Value = 4904.273
$ws.Range("I40").Value = 3394.75
$ws.Range("K40").Value = 3394.75
$ws.Range("M40").Value = -3258.75

$ws.Range("H46").Value = 12266.833
$ws.Range("J46").Value = 18684.5
$ws.Range("L46").Value = 18684.5
$ws.Range("N46").Value = -19060.5

$ws.Range("H93").Value = 1654.1538
$ws.Range("I93").Value = 1614.909
$ws.Range("K93").Value = 1614.909
$ws.Range("M93").Value = -366.9090000000001

$ws.Range("H132").Value = 5561.037
$ws.Range("J132").Value = 6331.778
$ws.Range("L132").Value = 18995.334
$ws.Range("N132").Value = -24055.334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H12").Value = 25000
$ws.Range("I12").Value = 25000
$ws.Range("K12").Value = 25000
$ws.Range("M12").Value = -24858

$ws.Range("H113").Value = 3969421
$ws.Range("I113").Value = 7576604.5
$ws.Range("J113").Value = 1518.7
$ws.Range("K113").Value = 22729813.5
$ws.Range("L113").Value = 4556.1
$ws.Range("M113").Value = -22727643.5
$ws.Range("N113").Value = -8896.1

$ws.Range("H132").Value = 2099.2856
$ws.Range("I132").Value = 2049.5
$ws.Range("K132").Value = 6148.5
$ws.Range("M132").Value = -3618.5

$ws.Range("H136").Value = 2417.5454
$ws.Range("I136").Value = 2099.6072
$ws.Range("K136").Value = 6298.821599999999
$ws.Range("M136").Value = -3748.821599999999
